$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.376852989196777
$ws.Range("B1").Value = 1.654488086700439
$ws.Range("C1").Value = 2.225682973861694
$ws.Range("D1").Value = 4.998605728149414
$ws.Range("E1").Value = 2.189740896224976
